function Set-CellText($ws, $cellRef, $text) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "28.130.28"
Set-CellText $ws "E2" "  -3.25%  "
Set-CellText $ws "D3" "1.928.34"
Set-CellText $ws "E3" "  -2.22%  "
Set-CellText $ws "D4" "1.002"
Set-CellText $ws "E4" "  -0.64%  "
Set-CellText $ws "D5" "329.03"
Set-CellText $ws "E5" "  +0.03%  "
Set-CellText $ws "E6" "  -0.53%  "
Set-CellText $ws "D7" "0.4732"
Set-CellText $ws "E7" "  -5.08%  "
Set-CellText $ws "D8" "0.4072"
Set-CellText $ws "E8" "  -3.41%  "
Set-CellText $ws "D9" "52.93"
Set-CellText $ws "E9" "  -0.54%  "
Set-CellText $ws "D10" "0.08458"
Set-CellText $ws "E10" "  -8.47%  "
Set-CellText $ws "D11" "1.049"
Set-CellText $ws "E11" "  -4.64%  "
Set-CellText $ws "D12" "22.32"
Set-CellText $ws "E12" "  -2.34%  "
Set-CellText $ws "D13" "1.933.95"
Set-CellText $ws "E13" "  -2.20%  "
Set-CellText $ws "D14" "7.528"
Set-CellText $ws "E14" "  -4.88%  "
Set-CellText $ws "D15" "6.106"
Set-CellText $ws "E15" "  -5.42%  "
Set-CellText $ws "D16" "1.002"
Set-CellText $ws "E16" "  -0.83%  "
Set-CellText $ws "D17" "90.53"
Set-CellText $ws "D18" "0.00001067"
Set-CellText $ws "E18" "  -3.63%  "
Set-CellText $ws "D19" "0.06584"
Set-CellText $ws "E19" "  -2.10%  "
Set-CellText $ws "D20" "18.15"
Set-CellText $ws "E20" "  -5.80%  "
Set-CellText $ws "E21" "  -0.36%  "
Set-CellText $ws "D22" "5.764"
Set-CellText $ws "E22" "  -3.35%  "
Set-CellText $ws "D23" "28.158.61"
Set-CellText $ws "E23" "  -3.30%  "
Set-CellText $ws "D24" "11.43"
Set-CellText $ws "E24" "  -4.45%  "
Set-CellText $ws "D25" "2.275"
Set-CellText $ws "E25" "  +0.62%  "
Set-CellText $ws "D26" "2.197.37"
Set-CellText $ws "E26" "  -0.66%  "
Set-CellText $ws "D27" "154.39"
Set-CellText $ws "E27" "  -0.55%  "
Set-CellText $ws "D28" "20.14"
Set-CellText $ws "E28" "  -2.93%  "
Set-CellText $ws "D29" "2.160"
Set-CellText $ws "E29" "  -4.64%  "
Set-CellText $ws "D30" "5.725"
Set-CellText $ws "E30" "  -9.70%  "
Set-CellText $ws "D31" "123.82"
Set-CellText $ws "E31" "  -2.44%  "
Set-CellText $ws "D32" "0.9805"
Set-CellText $ws "E32" "  -6.61%  "
Set-CellText $ws "D33" "0.09613"
Set-CellText $ws "E33" "  -2.51%  "
Set-CellText $ws "D34" "1.452"
Set-CellText $ws "E34" "  -4.51%  "
Set-CellText $ws "D35" "5.570"
Set-CellText $ws "E35" "  -4.23%  "
Set-CellText $ws "D36" "3.637"
Set-CellText $ws "E36" "  -2.54%  "
Set-CellText $ws "D37" "9.098"
Set-CellText $ws "E37" "  +0.49%  "
Set-CellText $ws "D38" "0.02325"
Set-CellText $ws "E38" "  -4.37%  "
Set-CellText $ws "D39" "0.06183"
Set-CellText $ws "E39" "  -3.80%  "
Set-CellText $ws "D40" "1.241"
Set-CellText $ws "E40" "  -6.13%  "
Set-CellText $ws "D41" "0.6185"
Set-CellText $ws "E41" "  -4.61%  "
Set-CellText $ws "D42" "11.09"
Set-CellText $ws "E42" "  -3.71%  "
Set-CellText $ws "D43" "1.002"
Set-CellText $ws "E43" "  -0.52%  "
Set-CellText $ws "D44" "0.1909"
Set-CellText $ws "E44" "  -4.64%  "
Set-CellText $ws "D45" "0.5909"
Set-CellText $ws "E45" "  -5.12%  "
Set-CellText $ws "D46" "1.296"
Set-CellText $ws "E46" "  -5.45%  "
Set-CellText $ws "D47" "12.79"
Set-CellText $ws "E47" "  -3.86%  "
Set-CellText $ws "D48" "2.040"
Set-CellText $ws "E48" "  -6.80%  "
Set-CellText $ws "E49" "  -0.29%  "
Set-CellText $ws "D50" "0.06820"
Set-CellText $ws "E50" "  -2.01%  "
Set-CellText $ws "D51" "110.19"
Set-CellText $ws "E51" "  -2.76%  "
